$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Delete the duplicate "Contact" row (row 11) entirely, shifting rows below up by one.
$ws.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: update timestamp
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$ws.Range("B9").Value = "Alvearie Team"

# Former "Contact" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value (force plain text "true", not a boolean)
$ws.Range("B14").Value = "'true"
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
